# weight_change.xlsx — "Add files via upload" re-edit
#
# 1) Shorten the header in C1 from " Physical Characteristics" to "Phys".
# 2) Update the saved view state: scroll the window right one column
#    (so column B becomes the left-most visible column) and move the
#    active selection from C30 to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edit -----------------------------------------------------
$ws.Range("C1").Value = "Phys"

# --- View / selection state --------------------------------------------
# Scroll so column B is the left-most visible column (topLeftCell = B1).
$excel.ActiveWindow.ScrollColumn = 2

# Move the active cell / selection to C8.
$ws.Range("C8").Select()
